$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a total column (G) that sums the covariate counts in columns C:F for row 4
$ws.Range("G4").Formula = "=SUM(C4:F4)"

# Match the author's final selection, which moved to the new total cell
$ws.Range("G4").Select() | Out-Null
